# "Repayment schedule" sheet: insert a new (blank) column before column N
# ("Late"), pushing Late / heading / Outstanding one column to the right,
# and give the new column the same width as the column to its left (M) -
# exactly what Excel does on a manual column insert.
$wb = $excel.ActiveWorkbook
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

$wsRepay.Activate() | Out-Null

$wsRepay.Columns("N").Insert()
$wsRepay.Columns("N").ColumnWidth = $wsRepay.Columns("M").ColumnWidth

# Move the selection/active cell on the Repayment schedule sheet, which
# also makes it the active (selected) tab of the workbook instead of
# "Transactions".
$wsRepay.Range("R8").Select() | Out-Null
